# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Terminal Hortofrutícola Agro
# Chillán - Frutilla" at row 104 (pushing the existing rows 104:120 down
# to 105:121), then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 104; this shifts rows 104-120
# down to 105-121 and extends the sheet dimension to A1:T121.
$ws.Rows.Item(104).Insert()

$ws.Range("A104").Value = 7
$ws.Range("B104").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C104").Value = "Ñuble"
$ws.Range("D104").Value = 44476
$ws.Range("E104").Value = 16
$ws.Range("F104").Value = "Fruta"
$ws.Range("G104").Value = 100101
$ws.Range("H104").Value = "Berries"
$ws.Range("I104").Value = 100112025
$ws.Range("J104").Value = "Frutilla"
$ws.Range("K104").Value = "Sin especificar"
$ws.Range("L104").Value = "Segunda"
$ws.Range("M104").Value = 60
$ws.Range("N104").Value = 12000
$ws.Range("O104").Value = 13000
$ws.Range("P104").Value = 12500
$ws.Range("Q104").Value = "`$/bandeja 7 kilos"
$ws.Range("R104").Value = "Provincia de Melipilla"
$ws.Range("S104").Value = 1786
$ws.Range("T104").Value = 7
